$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cells I1 and J1 with same formatting as the existing header H1
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data for rows 2..57: columns I (I0) and J (IF)
$data = @(
    @(3,3),
    @(6,8),
    @(6,7),
    @(8,8),
    @(8,8),
    @(6,6),
    @(8,8),
    @(7,7),
    @(7,8),
    @(7,7),
    @(6,7),
    @(8,8),
    @(7,7),
    @(8,8),
    @(6,6),
    @(7,7),
    @(7,7),
    @(6,6),
    @(6,6),
    @(5,6),
    @(7,8),
    @(6,7),
    @(7,7),
    @(10,10),
    @(6,6),
    @(7,8),
    @(6,6),
    @(7,7),
    @(7,8),
    @(7,7),
    @(7,8),
    @(6,6),
    @(8,8),
    @(7,7),
    @(6,7),
    @(7,7),
    @(7,7),
    @(8,8),
    @(7,7),
    @(7,8),
    @(7,7),
    @(7,7),
    @(7,7),
    @(8,8),
    @(8,8),
    @(5,6),
    @(8,8),
    @(6,6),
    @(7,7),
    @(6,6),
    @(7,7),
    @(6,6),
    @(6,6),
    @(8,8),
    @(6,6),
    @(3,3)
)

for ($k = 0; $k -lt $data.Length; $k++) {
    $row = $k + 2
    $ws.Cells.Item($row, 9).Value = $data[$k][0]
    $ws.Cells.Item($row, 10).Value = $data[$k][1]
}
